$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "We are seeking a Software Engineer to build and maintain high-quality software solutions." + [char]10 + "Work with global teams to drive innovation and deliver scalable applications." + [char]10 + "Join Akkodis and be part of a tech-driven, collaborative environment."

$ws.Cells.Item(7,1).Value = "JD_006"
$ws.Cells.Item(7,2).Value = "Senior IT Engineer"
$ws.Cells.Item(7,3).Value = $desc
$ws.Cells.Item(7,4).Value = 5
$ws.Cells.Item(7,5).Value = 7

$ws.Rows.Item(7).AutoFit()
